{"js": "// Replace the 25 three-digit-by-one-digit multiplication equations found in\n// the single table of the document. Each cell's old equation text is\n// replaced with its new equation text while preserving all existing run\n// formatting (font, size, paragraph alignment, etc.) by using a targeted\n// search-and-replace (\"Replace\" insert location) scoped to each individual\n// table cell. Scoping to the cell (rather than searching the whole body)\n// is required because some equation strings (e.g. \"602\u00d73=1806\") occur more\n// than once in the document but must map to different replacement values\n// depending on which cell they are in.\n\nconst replacements = [\n  [0, 0, \"171\u00d72=342\", \"801\u00d72=1602\"],\n  [0, 1, \"151\u00d78=1208\", \"458\u00d76=2748\"],\n  [0, 2, \"981\u00d73=2943\", \"559\u00d73=1677\"],\n  [0, 3, \"678\u00d79=6102\", \"888\u00d74=3552\"],\n  [0, 4, \"245\u00d77=1715\", \"567\u00d75=2835\"],\n  [4, 0, \"465\u00d75=2325\", \"191\u00d72=382\"],\n  [4, 1, \"848\u00d76=5088\", \"112\u00d73=336\"],\n  [4, 2, \"602\u00d73=1806\", \"252\u00d76=1512\"],\n  [4, 3, \"581\u00d74=2324\", \"377\u00d79=3393\"],\n  [4, 4, \"365\u00d77=2555\", \"138\u00d73=414\"],\n  [9, 0, \"709\u00d74=2836\", \"227\u00d73=681\"],\n  [9, 1, \"299\u00d76=1794\", \"520\u00d75=2600\"],\n  [9, 2, \"525\u00d79=4725\", \"416\u00d73=1248\"],\n  [9, 3, \"340\u00d77=2380\", \"279\u00d73=837\"],\n  [9, 4, \"699\u00d78=5592\", \"835\u00d77=5845\"],\n  [14, 0, \"554\u00d75=2770\", \"480\u00d74=1920\"],\n  [14, 1, \"507\u00d78=4056\", \"233\u00d78=1864\"],\n  [14, 2, \"693\u00d79=6237\", \"485\u00d77=3395\"],\n  [14, 3, \"602\u00d73=1806\", \"261\u00d75=1305\"],\n  [14, 4, \"443\u00d79=3987\", \"617\u00d73=1851\"],\n  [19, 0, \"837\u00d74=3348\", \"672\u00d74=2688\"],\n  [19, 1, \"720\u00d78=5760\", \"431\u00d76=2586\"],\n  [19, 2, \"788\u00d73=2364\", \"226\u00d76=1356\"],\n  [19, 3, \"414\u00d73=1242\", \"410\u00d75=2050\"],\n  [19, 4, \"231\u00d74=924\", \"277\u00d77=1939\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const [row, col, oldText, newText] of replacements) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 three-digit-by-one-digit multiplication equations found in\n# the single table of the document. Each cell's old equation text is\n# replaced with its new equation text while preserving all existing run\n# formatting (font, size, paragraph alignment, etc.) by using Find/Replace\n# scoped to each individual table cell's Range. Scoping to the cell (rather\n# than running Find/Replace over the whole document) is required because\n# some equation strings (e.g. \"602\u00d73=1806\") occur more than once in the\n# document but must map to different replacement values depending on which\n# cell they are in.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$replacements = @(\n    @(1, 1, \"171\u00d72=342\", \"801\u00d72=1602\"),\n    @(1, 2, \"151\u00d78=1208\", \"458\u00d76=2748\"),\n    @(1, 3, \"981\u00d73=2943\", \"559\u00d73=1677\"),\n    @(1, 4, \"678\u00d79=6102\", \"888\u00d74=3552\"),\n    @(1, 5, \"245\u00d77=1715\", \"567\u00d75=2835\"),\n    @(5, 1, \"465\u00d75=2325\", \"191\u00d72=382\"),\n    @(5, 2, \"848\u00d76=5088\", \"112\u00d73=336\"),\n    @(5, 3, \"602\u00d73=1806\", \"252\u00d76=1512\"),\n    @(5, 4, \"581\u00d74=2324\", \"377\u00d79=3393\"),\n    @(5, 5, \"365\u00d77=2555\", \"138\u00d73=414\"),\n    @(10, 1, \"709\u00d74=2836\", \"227\u00d73=681\"),\n    @(10, 2, \"299\u00d76=1794\", \"520\u00d75=2600\"),\n    @(10, 3, \"525\u00d79=4725\", \"416\u00d73=1248\"),\n    @(10, 4, \"340\u00d77=2380\", \"279\u00d73=837\"),\n    @(10, 5, \"699\u00d78=5592\", \"835\u00d77=5845\"),\n    @(15, 1, \"554\u00d75=2770\", \"480\u00d74=1920\"),\n    @(15, 2, \"507\u00d78=4056\", \"233\u00d78=1864\"),\n    @(15, 3, \"693\u00d79=6237\", \"485\u00d77=3395\"),\n    @(15, 4, \"602\u00d73=1806\", \"261\u00d75=1305\"),\n    @(15, 5, \"443\u00d79=3987\", \"617\u00d73=1851\"),\n    @(20, 1, \"837\u00d74=3348\", \"672\u00d74=2688\"),\n    @(20, 2, \"720\u00d78=5760\", \"431\u00d76=2586\"),\n    @(20, 3, \"788\u00d73=2364\", \"226\u00d76=1356\"),\n    @(20, 4, \"414\u00d73=1242\", \"410\u00d75=2050\"),\n    @(20, 5, \"231\u00d74=924\", \"277\u00d77=1939\")\n)\n\nforeach ($item in $replacements) {\n    $rowIndex = $item[0]\n    $colIndex = $item[1]\n    $oldText = $item[2]\n    $newText = $item[3]\n\n    $cell = $tbl.Cell($rowIndex, $colIndex)\n    $range = $cell.Range\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $oldText\n    $range.Find.Replacement.Text = $newText\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n}\n"}
